# Weekly update: shift the existing Pomelo price records down by 3 rows
# (making room for 3 new weekly observations) by inserting 3 new rows
# right before the current row 248 equivalent (i.e. before the block
# shifts, the insertion point is at row 245, which is where the data
# currently starts).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("245:247").Insert()

# --- New row 245 ---
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 44776
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100102
$ws.Range("H245").Value = "Cítricos"
$ws.Range("I245").Value = 100102006
$ws.Range("J245").Value = "Pomelo"
$ws.Range("K245").Value = "Start Ruby"
$ws.Range("L245").Value = "Especial"
$ws.Range("M245").Value = 30
$ws.Range("N245").Value = 15000
$ws.Range("O245").Value = 15000
$ws.Range("P245").Value = 15000
$ws.Range("Q245").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R245").Value = "Región de O'Higgins"
$ws.Range("S245").Value = 1000
$ws.Range("T245").Value = 15

# --- New row 246 ---
$ws.Range("A246").Value = 10
$ws.Range("B246").Value = "Vega Modelo de Temuco"
$ws.Range("C246").Value = "La Araucanía"
$ws.Range("D246").Value = 44776
$ws.Range("E246").Value = 9
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100102
$ws.Range("H246").Value = "Cítricos"
$ws.Range("I246").Value = 100102006
$ws.Range("J246").Value = "Pomelo"
$ws.Range("K246").Value = "Start Ruby"
$ws.Range("L246").Value = "Primera"
$ws.Range("M246").Value = 100
$ws.Range("N246").Value = 10000
$ws.Range("O246").Value = 12000
$ws.Range("P246").Value = 11000
$ws.Range("Q246").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R246").Value = "Región de O'Higgins"
$ws.Range("S246").Value = 733
$ws.Range("T246").Value = 15

# --- New row 247 ---
$ws.Range("A247").Value = 10
$ws.Range("B247").Value = "Vega Modelo de Temuco"
$ws.Range("C247").Value = "La Araucanía"
$ws.Range("D247").Value = 44776
$ws.Range("E247").Value = 9
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100102
$ws.Range("H247").Value = "Cítricos"
$ws.Range("I247").Value = 100102006
$ws.Range("J247").Value = "Pomelo"
$ws.Range("K247").Value = "Start Ruby"
$ws.Range("L247").Value = "Segunda"
$ws.Range("M247").Value = 30
$ws.Range("N247").Value = 8000
$ws.Range("O247").Value = 8000
$ws.Range("P247").Value = 8000
$ws.Range("Q247").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R247").Value = "Región de O'Higgins"
$ws.Range("S247").Value = 533
$ws.Range("T247").Value = 15
